# Regenerate the localization handoff report:
#  - Priority for the "Ready for handoff" rows (4-7) bumps from "low" to "ht"
#  - The "Latest Handoff Datetime" for those same rows is refreshed to the
#    new report-generation timestamps (one value per language sheet).
$wb = $excel.ActiveWorkbook

$rows = @(4, 5, 6, 7)

$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "ht"
    $ws.Range("H$r").Value = "2016-08-22 10:11:21"
}

$ws = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "ht"
    $ws.Range("H$r").Value = "2016-08-22 10:11:26"
}

# The Overview sheet's "Latest HO Xliff Generate Date" mirrors de-de's
# "Latest Handoff Datetime" for these rows, so refresh it too.
$ws = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "2016-08-22 10:11:26"
}
